$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64 - this shifts rows 64:118 down to 65:119
# (the previous last row, old row 118, becomes the new row 119).
$ws.Rows("64:64").Insert()

# Populate the freshly inserted row 64 with the new data point.
$ws.Range("A64").Value = 10
$ws.Range("B64").Value = "Vega Modelo de Temuco"
$ws.Range("C64").Value = "La Araucanía"
$ws.Range("D64").Value = 45096
$ws.Range("E64").Value = 9
$ws.Range("F64").Value = 100112010
$ws.Range("G64").Value = "Achicoria"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 55
$ws.Range("K64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("M64").Value = 9000
$ws.Range("N64").Value = "$/caja 18 unidades"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 500
$ws.Range("Q64").Value = 18
$ws.Range("R64").Value = "Hortaliza"
